# Standardize Processing Question Display -1 10 03 25
# Append new Q&A rows (184-197) to Sheet1, matching original inline-string formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('What the maximum number of headers I can display in my log?', 'llama3.2:latest', 'The maximum number of headers you can display in your log is up to 50.'),
    @('How many tables can I have in my log?', 'llama3.2:latest', 'You can have up to 100 tables in a log.'),
    @('How many symbols can I have in the plot at any one time?', 'llama3.2:latest', 'You can have up to 10,000 symbols in a plot at any one time.'),
    @('What is the maximum number of headers I can display in my log?', 'llama3.2:latest', 'The maximum number of headers you can display in your log is up to 50.'),
    @('What''s the maximum number of lithology types in an log?', 'llama3.2:latest', 'The highest number of lithology types that can be represented in a log is 450.'),
    @('How many tracks can you define in one ODF?', 'llama3.2:latest', 'The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most.'),
    @('How many curve shades can I create?', 'llama3.2:latest', 'You can create up to 250 curve shades in one ODF.'),
    @('How many curves can I load in one go?', 'llama3.2:latest', 'According to the document, you can load 450 curves in one go.'),
    @('What the maximum number of headers I can display in my log?', 'llama3.2:latest', 'The maximum number of headers you can display in your log is up to 50.'),
    @('How many tables can I have in my log?', 'llama3.2:latest', 'You can have up to 100 tables in a log.'),
    @('Whats the maximum number of characters in a single text entry?', 'llama3.2:latest', '32000.'),
    @('How many symbols can I have in the plot at any one time?', 'llama3.2:latest', 'You can have up to 10,000 symbols in a plot at any one time.'),
    @('How many scales can I define?', 'llama3.2:latest', 'According to the document, you can define up to 23 scales.'),
    @('What the maximum number of data files I can load?', 'llama3.2:latest', 'The maximum number of data files you can load is unlimited.'),
)

$startRow = 184
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    $cellC = $ws.Cells.Item($r, 3)
    $text = $row[2]
    # Force purely-numeric-looking strings (e.g. "32000.") to stay text,
    # then strip the temporary formatting so no style index is left on the cell.
    if ($text -match "^[0-9.,]+$") {
        $cellC.NumberFormat = "@"
        $cellC.Value = $text
        $cellC.ClearFormats()
    } else {
        $cellC.Value = $text
    }
}
